$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "42.463.20"
$ws.Range("E2").Value2 = "  +1.25%  "

$ws.Range("D3").Value2 = "2.305.91"
$ws.Range("E3").Value2 = "  +0.49%  "

$ws.Range("E4").Value2 = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "318.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +2.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "104.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +0.26%  "

$ws.Range("E7").Value2 = "  +1.40%  "

$ws.Range("E8").Value2 = "  -0.06%  "

$ws.Range("E9").Value2 = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "40.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +0.79%  "

$ws.Range("E11").Value2 = "  +0.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "8.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +1.89%  "

$ws.Range("E13").Value2 = "  +1.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.971"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "15.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +0.03%  "

$ws.Range("D16").Value2 = "2.655.15"
$ws.Range("E16").Value2 = "  +0.71%  "

$ws.Range("D17").Value2 = "2.306.82"
$ws.Range("E17").Value2 = "  +0.00%  "

$ws.Range("D18").Value2 = "42.607.49"
$ws.Range("E18").Value2 = "  +1.62%  "

$ws.Range("E19").Value2 = "  -1.11%  "

$ws.Range("E20").Value2 = "  +2.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "73.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "280.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +9.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "3.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +3.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "11.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +21.09%  "

$ws.Range("E25").Value2 = "  +1.74%  "

$ws.Range("E26").Value2 = "  -0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -0.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +6.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "22.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +1.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "36.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +2.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "165.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +1.26%  "

$ws.Range("E32").Value2 = "  -0.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "5.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +2.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.137"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +6.79%  "

$ws.Range("B35").Value2 = "Kaspa"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.117"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +1.74%  "

$ws.Range("B36").Value2 = "WEMIXToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -10.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.0376"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +7.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "4.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +3.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "3.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +3.58%  "

$ws.Range("E40").Value2 = "  +3.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +4.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "98.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "70.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +0.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -0.28%  "

$ws.Range("E45").Value2 = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "12.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "80.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +8.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "112.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +1.49%  "

$ws.Range("E49").Value2 = "  +0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "5.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -0.56%  "

$ws.Range("D51").Value2 = "1.604.61"
$ws.Range("E51").Value2 = "  +4.71%  "
